$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A9").Value = 104
$ws.Range("B9").Value = "LC/GFG"
$ws.Range("C9").Value = "Maximum Depth of Binary Tree or Height of Binary Tree"

$ws.Range("A7:C7").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C13").Select()
